$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# "split evaluating and predicting": the table headers already existed, but the
# actual evaluation numbers (F1 / BLEU scores) were never filled in for the four
# decoding methods. Backfill them now, each as a 2-decimal number.
$values = @(
    @(0.34182366120735902, 58.5817266068387,  54.371892507957298, 19.561407387726302, 19.511078188773201),
    @(0.28147622203288197, 57.2424796521845,  48.4727810559828,   18.353755153828999, 17.279835598492699),
    @(0.25959481207990098, 48.467653289467698, 49.078043567809303, 14.607165790417501, 14.5870804303157),
    @(0.283959733661522,   53.443229036504299, 52.246567667076,    15.509973904627,    15.538847615420099)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $r + 2
    $rowValues = $values[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $col = $c + 2
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $rowValues[$c]
        $cell.NumberFormat = "0.00"
    }
}

# Column F now holds real numeric data instead of being empty, so it needs to be
# widened to fit.
$ws.Columns.Item(6).ColumnWidth = 30

# Selection moved from H12 to D12
$ws.Range("D12").Select()
